$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 (gene FAM86A): "Unknown" -> "Psuedogene predicted to be involved in methylation"
$ws.Range("D34").Value = "Psuedogene predicted to be involved in methylation"

# Row 35 (gene MSI2): stays "Stem cell maintenance and cancer progression " (unchanged text)
$ws.Range("D35").Value = "Stem cell maintenance and cancer progression "

# Row 36 (gene CECR2): "Unknown" -> "Chromatin remodeling"
$ws.Range("D36").Value = "Chromatin remodeling"

# Update the active selection to D36
$ws.Range("D36").Select()
